$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 717
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 717
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H47").Value = 21000
$ws.Range("J47").Value = 21000
$ws.Range("L47").Value = 21000
$ws.Range("N47").Value = -22944

$ws.Range("H74").Value = 4206.3125
$ws.Range("I74").Value = 3829.9
$ws.Range("J74").Value = 4833.6665
$ws.Range("K74").Value = 3829.9
$ws.Range("L74").Value = 4833.6665
$ws.Range("M74").Value = -2893.9
$ws.Range("N74").Value = -6705.6665

$ws.Range("H76").Value = 3903.7273
$ws.Range("I76").Value = 3059.9
$ws.Range("J76").Value = 4606.9165
$ws.Range("K76").Value = 3059.9
$ws.Range("L76").Value = 4606.9165
$ws.Range("M76").Value = -2744.9
$ws.Range("N76").Value = -5236.9165

$ws.Range("H77").Value = 4206.3125
$ws.Range("I77").Value = 3829.9
$ws.Range("J77").Value = 4833.6665
$ws.Range("K77").Value = 19149.5
$ws.Range("L77").Value = 24168.3325
$ws.Range("M77").Value = -14469.5
$ws.Range("N77").Value = -33528.3325

$ws.Range("H79").Value = 3903.7273
$ws.Range("I79").Value = 3059.9
$ws.Range("J79").Value = 4606.9165
$ws.Range("K79").Value = 3059.9
$ws.Range("L79").Value = 4606.9165
$ws.Range("M79").Value = -1967.9
$ws.Range("N79").Value = -6790.9165

$ws.Range("H110").Value = 60000
$ws.Range("J110").Value = 60000
$ws.Range("L110").Value = 60000
$ws.Range("N110").Value = -68180

$ws.Range("H132").Value = 3650.4102
$ws.Range("I132").Value = 2917.6875
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 8753.0625
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -6223.0625
$ws.Range("N132").Value = -26060

$ws.Range("H137").Value = 13334317
$ws.Range("I137").Value = 812.6667
$ws.Range("J137").Value = 66668332
$ws.Range("K137").Value = 2438.0001
$ws.Range("L137").Value = 200004996
$ws.Range("M137").Value = 111.9998999999998
$ws.Range("N137").Value = -200010096

$ws.Range("H141").Value = 3076.182
$ws.Range("I141").Value = 1106.8125
$ws.Range("J141").Value = 4929.706
$ws.Range("K141").Value = 3320.4375
$ws.Range("L141").Value = 14789.118
$ws.Range("M141").Value = 1859.5625
$ws.Range("N141").Value = -25149.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6424.99
$ws.Range("I32").Value = 4478.946
$ws.Range("J32").Value = 11963.73
$ws.Range("K32").Value = 4478.946
$ws.Range("L32").Value = 11963.73
$ws.Range("M32").Value = -4191.946
$ws.Range("N32").Value = -12537.73

$ws.Range("H44").Value = 28686.625
$ws.Range("J44").Value = 28686.625
$ws.Range("L44").Value = 28686.625
$ws.Range("N44").Value = -29662.625

$ws.Range("H55").Value = 34206.5
$ws.Range("J55").Value = 34206.5
$ws.Range("L55").Value = 34206.5
$ws.Range("N55").Value = -34836.5

$ws.Range("H63").Value = 100000760
$ws.Range("I63").Value = 100000760
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 100000760
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -100000074
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 100000760
$ws.Range("I66").Value = 100000760
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 500003800
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -500000368
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4058.2068
$ws.Range("I105").Value = 2082
$ws.Range("K105").Value = 2082
$ws.Range("M105").Value = -335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7579480.5
$ws.Range("J31").Value = 55559830
$ws.Range("L31").Value = 55559830
$ws.Range("N31").Value = -55560420

$ws.Range("H34").Value = 7579480.5
$ws.Range("J34").Value = 55559830
$ws.Range("L34").Value = 55559830
$ws.Range("N34").Value = -55560234

$ws.Range("H99").Value = 2133.3333
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2400
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2400
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5396

$ws.Range("H122").Value = 2372.6667
$ws.Range("I122").Value = 2716.25
$ws.Range("J122").Value = 1685.5
$ws.Range("K122").Value = 8148.75
$ws.Range("L122").Value = 5056.5
$ws.Range("M122").Value = -5698.75
$ws.Range("N122").Value = -9956.5

$ws.Range("H126").Value = 2133.3333
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -12140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 9000
$ws.Range("J5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("N5").Value = -9224

$ws.Range("H49").Value = 20874.75
$ws.Range("J49").Value = 20874.75
$ws.Range("L49").Value = 20874.75
$ws.Range("N49").Value = -21242.75

$ws.Range("H80").Value = 9011576
$ws.Range("I80").Value = 14494989
$ws.Range("J80").Value = 3112.7856
$ws.Range("K80").Value = 14494989
$ws.Range("L80").Value = 3112.7856
$ws.Range("M80").Value = -14493991
$ws.Range("N80").Value = -5108.7856

$ws.Range("H83").Value = 9011576
$ws.Range("I83").Value = 14494989
$ws.Range("J83").Value = 3112.7856
$ws.Range("K83").Value = 72474945
$ws.Range("L83").Value = 15563.928
$ws.Range("M83").Value = -72469953
$ws.Range("N83").Value = -25547.928

$ws.Range("H113").Value = 42668.082
$ws.Range("I113").Value = 72083.71000000001
$ws.Range("J113").Value = 1486.2
$ws.Range("K113").Value = 72083.71000000001
$ws.Range("L113").Value = 1486.2
$ws.Range("M113").Value = -69913.71000000001
$ws.Range("N113").Value = -5826.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 19053.889
$ws.Range("I56").Value = 2400
$ws.Range("K56").Value = 2400
$ws.Range("M56").Value = -1709
